# Auto-generated edit script: update cryptos price/volume columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.834.52"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "3.495.43"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.64"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.55"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.36%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -1.29%  "
$ws.Range("E9").Value = "  +4.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.18"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.03%  "
$ws.Range("E11").Value = "  -0.94%  "
$ws.Range("D12").Value = "4.097.79"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("E13").Value = "  +0.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.30"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +4.32%  "
$ws.Range("D15").Value = "66.837.13"
$ws.Range("E15").Value = "  +0.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000178"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("D17").Value = "3.496.80"
$ws.Range("E17").Value = "  +0.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.27"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.21"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "393.32"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.93"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.21"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.53%  "
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("E24").Value = "  +0.28%  "
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.22"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.994"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.13"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -3.02%  "
$ws.Range("E30").Value = "  -2.03%  "
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "23.62"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.34"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.61"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "162.68"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  -1.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.90"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.85"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.63"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").Value = "2.841.77"
$ws.Range("E40").Value = "  +2.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.18"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0737"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.06"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.67"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.91%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.53"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0302"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "338.34"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "34.58"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.48%  "
$ws.Range("E49").Value = "  -1.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.42"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.841"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.14%  "
